$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in the "Battery capacity" label (A1) ---
$ws.Range("A1").Value = "Battery capacity"

# --- Update current-consumption inputs used in the estimate ---
# Alarm on with LED current (mA)
$ws.Range("B6").Value = 3.5
# Number of times alarm goes off per week
$ws.Range("B11").Value = 2

# --- Add reviewer comment on the updated "Alarm on with LED" measurement ---
$comment = $ws.Range("B6").AddComment("Andreas:" + [char]10 + "Measured for 600 Hz with 10% duty cycle.")

# --- Restore the selected cell as left by the editor ---
$ws.Range("D30").Select()
